# [AE-880] first update with latest files
# - bump the cached "datetimeFigureOut" field text (1/23/14 -> 9/6/14) on
#   the slide master and every slide layout's Date Placeholder
# - relabel the streaming-figures diagram on slide 1: HDFS -> HDFS/S3,
#   ZeroMQ -> Kinesis

$p = $ppt.ActivePresentation

# ---- 1. Date placeholders on the master + all custom (slide) layouts ----
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "1/23/14") {
            $sh.TextFrame.TextRange.Text = "9/6/14"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# ---- 2. HDFS -> HDFS/S3 and ZeroMQ -> Kinesis on slide 1 ----
$slide1 = $p.Slides.Item(1)
$diagram = $slide1.Shapes.Item(1)

for ($i = 1; $i -le $diagram.GroupItems.Count; $i++) {
    $sh = $diagram.GroupItems.Item($i)
    if (-not $sh.HasTextFrame) { continue }
    $tr = $sh.TextFrame.TextRange
    if ($sh.Name -eq "Rounded Rectangle 47" -and $tr.Text -eq "HDFS") {
        $tr.Text = "HDFS/S3"
    } elseif ($sh.Name -eq "Rounded Rectangle 48" -and $tr.Text -eq "ZeroMQ") {
        $tr.Text = "Kinesis"
    }
}
